$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "#! END_ROW" / "#! END_ROW true" markers added in column J, one per
# existing row of the pyramid (rows 1-6). Write the plain "#! END_ROW" cells
# first so the shared-string table gets "#! END_ROW" before "#! END_ROW true".
$ws.Range("J2").Value = "#! END_ROW"
$ws.Range("J4").Value = "#! END_ROW"
$ws.Range("J6").Value = "#! END_ROW"

$ws.Range("J1").Value = "#! END_ROW true"
$ws.Range("J3").Value = "#! END_ROW true"
$ws.Range("J5").Value = "#! END_ROW true"

# Move the active selection like the author's session ended up (H19).
$null = $ws.Range("H19").Select()
